$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.936.69"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.291.62"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'252.35"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "'0.639"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("D7").Value = "'74.98"
$ws.Range("E7").Value = "  +7.00%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("D10").Value = "'39.26"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "'7.52"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").Value = "2.635.91"
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("D15").Value = "'15.08"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").Value = "2.290.45"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "42.833.82"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "'72.39"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "'236.01"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +5.11%  "
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'11.33"
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -7.20%  "
$ws.Range("D29").Value = "'167.45"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'21.06"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("E31").Value = "  +9.06%  "
$ws.Range("D32").Value = "'6.25"
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "'0.129"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("D39").Value = "'13.64"
$ws.Range("E39").Value = "  +9.25%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E42").Value = "  +4.24%  "
$ws.Range("D43").Value = "'9.11"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "'61.35"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'105.28"
$ws.Range("E46").Value = "  +11.03%  "
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "'4.24"
$ws.Range("E51").Value = "  -1.74%  "
